$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

$ws.Range("S2").Value = "'true"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 1

$ws.Range("S3").Value = "'true"
$ws.Range("S3").Style = "Normal"
$ws.Range("T3").Value = 1
